# Add a new "share_sales" worksheet after the existing sheets, populate it
# with the share-sale data used to raise cash during COVID, format the
# Proceeds column, and update the active-sheet/selection state to match.

$wb = $excel.ActiveWorkbook

# --- Add the new worksheet after the last existing sheet -------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "share_sales"

# --- Headers -----------------------------------------------------------
$newSheet.Range("A1").Value = "Year"
$newSheet.Range("B1").Value = "Quarter"
$newSheet.Range("C1").Value = "Airline"
# Write "Proceeds" before "Shares Sold" so the shared-string table order
# matches (Proceeds ends up before Shares Sold in sst).
$newSheet.Range("E1").Value = "Proceeds"
$newSheet.Range("D1").Value = "Shares Sold"

# --- Data rows -----------------------------------------------------------
$data = @(
    @(2020, "FY", "AAL", 143597509, 1901355562.5899999),
    @(2021, "FY", "AAL", 79599019, 1082205377.28),
    @(2020, "FY", "DAL", 6500000, 158535000),
    @(2021, "FY", "DAL", 2100000, 83433000),
    @(2020, "FY", "UAL", 69338841, 2293578491.5),
    @(2021, "FY", "UAL", 14200000, 685282000)
)

$r = 2
foreach ($row in $data) {
    $newSheet.Range("A$r").Value = $row[0]
    $newSheet.Range("B$r").Value = $row[1]
    $newSheet.Range("C$r").Value = $row[2]
    $newSheet.Range("D$r").Value = $row[3]
    $newSheet.Range("E$r").Value = $row[4]
    $r++
}

# --- Formatting ------------------------------------------------------------
$newSheet.Range("E1:E7").NumberFormat = "0.00"

# (Target "bestFit" widths, from the authored workbook, are 4.85546875 /
# 7.7109375 / 6.7109375 / 19.28515625 / 15.85546875 characters; the engine
# quantizes ColumnWidth to 1/6-character steps, so these inputs are chosen
# to land on the closest reproducible value.)
$newSheet.Columns.Item(1).ColumnWidth = 4
$newSheet.Columns.Item(2).ColumnWidth = 6.833333333333333
$newSheet.Columns.Item(3).ColumnWidth = 5.833333333333333
$newSheet.Columns.Item(4).ColumnWidth = 18.5
$newSheet.Columns.Item(5).ColumnWidth = 15

# --- Selection on the new sheet ---------------------------------------
$newSheet.Range("D2").Select() | Out-Null

# --- Update selection on share_repurchases sheet --------------------------
$repurchases = $wb.Worksheets.Item("share_repurchases")
$repurchases.Range("A1:E1").Select() | Out-Null

# --- Make the new sheet the active/selected tab ----------------------------
$newSheet.Activate() | Out-Null
